$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.873.68'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.82%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.810.70'
$ws.Range('D3').Style = "Normal"

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '309.43'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.32%  '

# Row 6
$ws.Range('E6').Value = '  +0.10%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4641'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.54%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3701'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.41%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07364'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.83%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8755'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.49%  '

# Row 11
$ws.Range('E11').Value = '  -1.90%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.846.74'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.46%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.364'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.54%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.509'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.31%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.07055'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.47%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.64'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.71%  '

# Row 17
$ws.Range('E17').Value = '  +0.08%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008704'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.40%  '

# Row 19
$ws.Range('E19').Value = '  +0.05%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.75'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.09%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '26.864.81'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.87%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.312'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.56%  '

# Row 23
$ws.Range('E23').Value = '  -3.94%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.963.47'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.06%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.898'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.79%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '151.57'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.29%  '

# Row 27
$ws.Range('E27').Value = '  -1.59%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.151'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.48%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.313'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.67%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '115.94'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.27%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08913'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.53%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7562'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.48%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.156'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.47%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.464'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.04%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.917'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.54%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.11%  '

# Row 37
$ws.Range('E37').Value = '  -0.19%  '

# Row 38
$ws.Range('E38').Value = '  -0.70%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.446'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.05%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.05262'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.11%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.911'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.65%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5319'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.56%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.169'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.29%  '

# Row 44
$ws.Range('E44').Value = '  -2.85%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.467'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.57%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4951'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.49%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.37'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.61%  '

# Row 48
$ws.Range('E48').Value = '  -0.90%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.001'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.13%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '103.00'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.46%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06290'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.54%  '
